# Incorporated review comments from Tanel

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Footer "datetimeFigureOut" field text: 8/17/2018 -> 8/27/2018
#    (slide master + every slide layout carries its own cached copy)
# ---------------------------------------------------------------------
$oldDate = "8/17/2018"
$newDate = "8/27/2018"

function Update-DateShape($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

$master = $p.SlideMaster
Update-DateShape $master.Shapes

for ($L = 1; $L -le $master.CustomLayouts.Count; $L++) {
    $layout = $master.CustomLayouts.Item($L)
    Update-DateShape $layout.Shapes
}

# ---------------------------------------------------------------------
# 2. Slide 2 title: "Outline" -> "Background and Topics Covered"
# ---------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$titleShape = $s2.Shapes.Item(1)
$titleShape.TextFrame.TextRange.Text = "Background and Topics Covered"

# ---------------------------------------------------------------------
# 3. Slide 8 content bullets
# ---------------------------------------------------------------------
$s8 = $p.Slides.Item(8)
$body = $s8.Shapes.Item(2)
$tr8 = $body.TextFrame.TextRange

# 3a. "Granularity of locking, reduces performance" ->
#     "Not having a fine grained granularity of locking, reduces performance"
$para3 = $tr8.Paragraphs(3, 1)
$para3.Text = "__TMP_PLACEHOLDER_PARA3__"
$para3 = $tr8.Paragraphs(3, 1)
$para3.Text = "Not having a fine grained granularity of locking, reduces performance"

# 3b. "Conditions for deadlock" ->
#     "Classical conditions for deadlock " + "to occur" (two runs)
$para5 = $tr8.Paragraphs(5, 1)
$para5.Text = "__TMP_PLACEHOLDER_PARA5__"
$para5 = $tr8.Paragraphs(5, 1)
$para5.Text = "Classical conditions for deadlock "
$para5 = $tr8.Paragraphs(5, 1)
$para5.InsertAfter("to occur") | Out-Null
